# bkiabi-skrining-preeklampsia.docx -- compatibility/formatting touch-up
#
# 1) The criteria table's two right-hand columns shrink/grow by a single
#    twip each (2203 -> 2202, 2117 -> 2118) across the grid definition and
#    every data cell in those columns.
# 2) The "( {nama_dokter} )" signature placeholder: the "{", "nama_dokter"
#    and "}" runs become one run, now bold, and the trailing ")" run has
#    its italic explicitly switched off.
#
# NOTE: paragraph indices (Paragraphs.Item(n)) go stale once the table
# geometry below is touched, so text is located purely via Range.Find
# instead of walking the Paragraphs collection.

$d = $word.ActiveDocument

# --- 1. Table column widths -------------------------------------------------
$table = $d.Tables.Item(1)
# Word reports/accepts column widths in points (1 pt = 20 dxa/twips).
$table.Columns.Item(2).Width = 2202 / 20   # 110.1 pt  (was 2203 dxa)
$table.Columns.Item(3).Width = 2118 / 20   # 105.9 pt  (was 2117 dxa)

# --- 2. Signature placeholder text/formatting -------------------------------

# Merge "{" + "nama_dokter" + "}" into a single bold run.
$find = $d.Content
$found = $find.Find.Execute("{nama_dokter}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $mergeStart = $find.Start
    $mergeLen = $find.End - $find.Start
    $find.Delete()
    $insertionPoint = $d.Range($mergeStart, $mergeStart)
    $insertionPoint.InsertAfter("{nama_dokter}")
    $mergedRun = $d.Range($mergeStart, $mergeStart + $mergeLen)
    $mergedRun.Font.Bold = $true

    # The trailing ")" (a few characters further on, in the same
    # paragraph) loses its italic formatting - now explicitly off.
    $afterMergeEnd = $mergeStart + $mergeLen
    $docEnd = $d.Content.End
    $searchEnd = [System.Math]::Min($afterMergeEnd + 10, $docEnd)
    $closeParen = $d.Range($afterMergeEnd, $searchEnd)
    $foundParen = $closeParen.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundParen) {
        $closeParen.Font.Italic = $false
    }
}
